# Updated cryptos list on Wed May 29 11:36:13 UTC 2024 with GitHub Actions
#
# This script updates the Price (column D) and Volume(1h) (column E) values
# for the crypto list in the active worksheet, and swaps the Mantle/Kaspa
# rows (40 and 41) per the latest scrape.
#
# All cells in columns D and E in this sheet are stored as plain text
# (many "prices" use '.' as a thousands separator, e.g. "67.880.17", and
# the percentages keep surrounding padding spaces). Excel's COM layer will
# silently reinterpret plain numeric-looking strings (e.g. "600.00",
# "18.70") as numbers and normalize/restyle them if we just set .Value.
# To avoid that, we force the cell to Text format before assignment and
# then restore the "Normal" style afterward so no stray formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Row, $Col, $Text) {
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextCell 2 4 "67.880.17"
Set-TextCell 2 5 "  -0.81%  "

# --- Row 3: Ethereum ---
Set-TextCell 3 4 "3.822.06"
Set-TextCell 3 5 "  -1.99%  "

# --- Row 4: TetherUSD ---
Set-TextCell 4 5 "  +0.01%  "

# --- Row 5: BNB ---
Set-TextCell 5 4 "600.00"
Set-TextCell 5 5 "  -0.38%  "

# --- Row 6: Solana ---
Set-TextCell 6 4 "169.08"
Set-TextCell 6 5 "  +0.08%  "

# --- Row 7: LidoStakedEther ---
Set-TextCell 7 4 "3.818.96"
Set-TextCell 7 5 "  -2.07%  "

# --- Row 8: USDC ---
Set-TextCell 8 5 "  +0.00%  "

# --- Row 9: XRP ---
Set-TextCell 9 4 "0.531"
Set-TextCell 9 5 "  +0.00%  "

# --- Row 10: Dogecoin ---
Set-TextCell 10 4 "0.165"
Set-TextCell 10 5 "  -0.55%  "

# --- Row 11: Toncoin ---
Set-TextCell 11 5 "  +1.19%  "

# --- Row 12: Cardano ---
Set-TextCell 12 5 "  +0.14%  "

# --- Row 13: ShibaInu ---
Set-TextCell 13 4 "0.0000275"

# --- Row 14: Avalanche ---
Set-TextCell 14 5 "  -0.39%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
Set-TextCell 15 4 "4.469.26"
Set-TextCell 15 5 "  -1.94%  "

# --- Row 16: WrappedEther ---
Set-TextCell 16 4 "3.811.18"
Set-TextCell 16 5 "  -2.47%  "

# --- Row 17: Chainlink ---
Set-TextCell 17 4 "18.70"
Set-TextCell 17 5 "  +2.59%  "

# --- Row 18: WrappedBTC ---
Set-TextCell 18 4 "68.036.20"
Set-TextCell 18 5 "  -0.65%  "

# --- Row 19: Polkadot ---
Set-TextCell 19 4 "7.46"
Set-TextCell 19 5 "  +0.25%  "

# --- Row 20: TRON ---
Set-TextCell 20 5 "  +0.18%  "

# --- Row 21: Uniswap ---
Set-TextCell 21 4 "10.86"
Set-TextCell 21 5 "  +0.11%  "

# --- Row 22: BitcoinCash ---
Set-TextCell 22 4 "468.86"
Set-TextCell 22 5 "  -0.76%  "

# --- Row 23: Polygon ---
Set-TextCell 23 5 "  -0.09%  "

# --- Row 24: PEPE ---
Set-TextCell 24 5 "  -9.13%  "

# --- Row 25: Litecoin ---
Set-TextCell 25 4 "83.77"
Set-TextCell 25 5 "  -0.07%  "

# --- Row 26: Fetch.AI ---
Set-TextCell 26 5 "  +2.27%  "

# --- Row 27: InternetComputer(DFINITY) ---
Set-TextCell 27 4 "12.23"
Set-TextCell 27 5 "  -0.01%  "

# --- Row 28: RenderToken ---
Set-TextCell 28 4 "10.43"
Set-TextCell 28 5 "  +4.14%  "

# --- Row 29: Dai ---
Set-TextCell 29 5 "  -0.09%  "

# --- Row 30: PancakeSwap ---
Set-TextCell 30 5 "  -1.42%  "

# --- Row 31: WrappedeETH ---
Set-TextCell 31 4 "3.974.76"
Set-TextCell 31 5 "  -1.93%  "

# --- Row 32: NEARProtocol ---
Set-TextCell 32 4 "7.77"
Set-TextCell 32 5 "  -1.80%  "

# --- Row 33: ImmutableX ---
Set-TextCell 33 5 "  -1.66%  "

# --- Row 34: EthereumClassic ---
Set-TextCell 34 4 "30.79"
Set-TextCell 34 5 "  -2.22%  "

# --- Row 35: Aptos ---
Set-TextCell 35 4 "9.33"
Set-TextCell 35 5 "  -1.57%  "

# --- Row 36: RenzoRestakedETH ---
Set-TextCell 36 4 "3.793.06"
Set-TextCell 36 5 "  -2.10%  "

# --- Row 37: dogwifhat ---
Set-TextCell 37 4 "3.88"
Set-TextCell 37 5 "  +4.46%  "

# --- Row 38: Hedera ---
Set-TextCell 38 5 "  +0.80%  "

# --- Row 39: Filecoin ---
Set-TextCell 39 4 "5.98"
Set-TextCell 39 5 "  +0.58%  "

# --- Rows 40 & 41: Mantle and Kaspa swap rank positions ---
# Row 40 was Mantle, becomes Kaspa; Row 41 was Kaspa, becomes Mantle.
Set-TextCell 40 2 "Kaspa"
Set-TextCell 40 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell 40 4 "0.139"
Set-TextCell 40 5 "  -1.65%  "

Set-TextCell 41 2 "Mantle"
Set-TextCell 41 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell 41 4 "1.01"
Set-TextCell 41 5 "  -1.75%  "

# --- Row 42: FirstDigitalUSD ---
Set-TextCell 42 4 "0.999"
Set-TextCell 42 5 "  -0.16%  "

# --- Row 43: TheGraph ---
Set-TextCell 43 5 "  +1.60%  "

# --- Row 44: USDe ---
Set-TextCell 44 5 "  -0.02%  "

# --- Row 45: Cosmos ---
Set-TextCell 45 4 "8.80"
Set-TextCell 45 5 "  +1.73%  "

# --- Row 46: Stacks ---
Set-TextCell 46 5 "  -0.87%  "

# --- Row 47: Bittensor ---
Set-TextCell 47 4 "409.05"
Set-TextCell 47 5 "  -4.45%  "

# --- Row 48: OKB ---
Set-TextCell 48 4 "46.58"
Set-TextCell 48 5 "  -1.25%  "

# --- Row 49: FLOKI ---
Set-TextCell 49 4 "0.000284"
Set-TextCell 49 5 "  -5.61%  "

# --- Row 50: Monero ---
Set-TextCell 50 4 "143.52"
Set-TextCell 50 5 "  -0.12%  "

# --- Row 51: VeChain ---
Set-TextCell 51 5 "  +0.06%  "
